$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 120, shifting existing rows 120:225 down to 121:226
$ws.Rows(120).EntireRow.Insert()

# Populate the new row 120 with its data (same fixed columns as every other
# data row in this sheet, plus the new record's own values).
$ws.Range("A120").Value = 9
$ws.Range("B120").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C120").Value = "Metropolitana"
$ws.Range("D120").Value = 44574
$ws.Range("E120").Value = 13
$ws.Range("F120").Value = 100112001
$ws.Range("G120").Value = "Berenjena"
$ws.Range("H120").Value = "Sin especificar"
$ws.Range("I120").Value = "Primera"
$ws.Range("J120").Value = 106
$ws.Range("K120").Value = 11000
$ws.Range("L120").Value = 13000
$ws.Range("M120").Value = 12000
$ws.Range("N120").Value = "$/caja 60 unidades"
$ws.Range("O120").Value = "Región Metropolitana"
$ws.Range("P120").Value = 200
$ws.Range("Q120").Value = 60
$ws.Range("R120").Value = "Hortaliza"

# Preserve the date number-format used by column D, matching the format
# already applied to every other date cell in column D.
$ws.Range("D120").NumberFormat = $ws.Range("D121").NumberFormat
